$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

foreach ($row in $table.Rows) {
    foreach ($cell in $row.Cells) {
        $rng = $cell.Range
        $rng.Find.ClearFormatting()
        $rng.Find.Execute("{", $true, $false, $false, $false, $false, $true, 1, $false, "<<", 2)
    }
}

foreach ($row in $table.Rows) {
    foreach ($cell in $row.Cells) {
        $rng = $cell.Range
        $rng.Find.ClearFormatting()
        $rng.Find.Execute("}", $true, $false, $false, $false, $false, $true, 1, $false, ">>", 2)
    }
}
